$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after the current row 98, pushing the old
# rows 99, 100, 101 down to 101, 102, 103 (their content stays intact).
$ws.Rows("99:100").Insert()

# Row 97: update in place with the new weekly reading.
$ws.Cells.Item(97, 4).Value = 45267
$ws.Cells.Item(97, 14).Value = 30000
$ws.Cells.Item(97, 15).Value = 30000
$ws.Cells.Item(97, 16).Value = 30000
$ws.Cells.Item(97, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(97, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(97, 20).Value = 15

# Row 98: update in place with the new weekly reading.
$ws.Cells.Item(98, 4).Value = 45267
$ws.Cells.Item(98, 11).Value = "Dina"
$ws.Cells.Item(98, 13).Value = 65
$ws.Cells.Item(98, 14).Value = 30000
$ws.Cells.Item(98, 15).Value = 30000
$ws.Cells.Item(98, 16).Value = 30000
$ws.Cells.Item(98, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(98, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(98, 19).Value = 2000
$ws.Cells.Item(98, 20).Value = 15

# Row 99 (newly inserted): fill with the data that previously lived in row 97.
$row99 = @(10, "Vega Modelo de Temuco", "La Araucanía", 44525, 9, "Fruta", 100103, "Frutos de hueso (carozo)", 100103003, "Damasco", "Castle Brite", "Primera", 55, 20000, 20000, 20000, "$/bandeja 10 kilos", "Provincia de San Felipe de Aconcagua", 2000, 10)
for ($i = 0; $i -lt $row99.Length; $i++) {
    $ws.Cells.Item(99, $i + 1).Value = $row99[$i]
}

# Row 100 (newly inserted): fill with the data that previously lived in row 98.
$row100 = @(10, "Vega Modelo de Temuco", "La Araucanía", 44559, 9, "Fruta", 100103, "Frutos de hueso (carozo)", 100103003, "Damasco", "Modesto", "Primera", 95, 18000, 18000, 18000, "$/bandeja 18 kilos", "Provincia de Quillota", 1000, 18)
for ($i = 0; $i -lt $row100.Length; $i++) {
    $ws.Cells.Item(100, $i + 1).Value = $row100[$i]
}
